# Errata table: insert a new row for "Percussion 2", bar 36, asking whether
# the final note of the bar is C or D. The row belongs right after the
# existing "Contrabassoon" / bar 35 row and right before the "Trumpets" /
# bar 37 row.

$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

# Locate the row that should come right *after* the new row (the
# "Trumpets", bar 37 row) by scanning column 1/2 text instead of hard
# coding a row index, so the script is resilient to any other table
# layout changes.
$targetRowIndex = -1
for ($i = 1; $i -le $t.Rows.Count; $i++) {
    $instrument = $t.Cell($i, 1).Range.Text.TrimEnd([char]13, [char]7)
    $bar = $t.Cell($i, 2).Range.Text.TrimEnd([char]13, [char]7)
    if ($instrument -eq "Trumpets" -and $bar -eq "37") {
        $targetRowIndex = $i
        break
    }
}

if ($targetRowIndex -eq -1) {
    throw "Could not locate the Trumpets/bar 37 row to insert before."
}

# Adding a row "before" an existing one copies that row's formatting
# (fonts, centering, cell widths, etc.), which matches the surrounding
# rows in this table.
$refRow = $t.Rows.Item($targetRowIndex)
$newRow = $t.Rows.Add($refRow)
$newRowIndex = $newRow.Index

$t.Cell($newRowIndex, 1).Range.Text = "Percussion 2"
$t.Cell($newRowIndex, 2).Range.Text = "36"
$t.Cell($newRowIndex, 3).Range.Text = "To clarify: is the final note in this bar C or D?"
# 4th column ("Answer") is intentionally left blank, matching the other
# unanswered rows in the table.
